$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Hora"
$ws.Range("C1").Value = "WC47 NACP"
$ws.Range("D1").Value = "WC48 POWER 5F"
$ws.Range("E1").Value = "WC49 POWER 5H"
$ws.Range("F1").Value = "WV50 FILTER"
$ws.Range("G1").Value = "SPL"

# Append new data rows 22-32
$data = @(
    @("2024-05-13", "12:10:35", "Fallo fijador tapa", "-", "-", "-", "-"),
    @("2024-05-13", "12:10:38", "-", "Etiquetadora", "-", "-", "-"),
    @("2024-05-13", "12:10:40", "-", "Etiquetadora", "-", "-", "-"),
    @("2024-05-13", "12:10:43", "-", "-", "-", "Traza", "-"),
    @("2024-05-13", "12:10:48", "Fallo cámara visión", "-", "-", "-", "-"),
    @("2024-05-13", "12:10:50", "AOI no detecta pieza", "-", "-", "-", "-"),
    @("2024-05-13", "12:10:52", "AOI no detecta pieza", "-", "-", "-", "-"),
    @("2024-05-13", "12:12:57", "AOI no detecta pieza", "-", "-", "-", "-"),
    @("2024-05-13", "12:13:01", "No coge placa", "-", "-", "-", "-"),
    @("2024-05-13", "12:13:07", "-", "Fallo etiqueta", "-", "-", "-"),
    @("2024-05-13", "12:13:10", "-", "-", "-", "Cover atascado", "-")
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    # Column A holds a date-looking string ("2024-05-13"); prefix with an
    # apostrophe so Excel stores it as literal text instead of converting it
    # to a date serial number, then reset the style so no extra
    # quote-prefix formatting sticks to the cell.
    $ws.Cells.Item($row, 1).Value = "'" + $data[$i][0]
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
    $ws.Cells.Item($row, 6).Value = $data[$i][5]
    $ws.Cells.Item($row, 7).Value = $data[$i][6]
}
